# Add a new "Padding in CNN" quiz question + answer slide pair, following
# the same authoring pattern already used throughout this deck: slide 27
# was an empty Title+Content "template" slide reserved for the next quiz
# question. We duplicate it twice (so the deck keeps a fresh empty
# template at the end, exactly like the other slide in the list), turn
# the original into the question slide and the first duplicate into the
# answer slide, dropping the (unused) Title placeholder from each.

$p = $ppt.ActivePresentation

$s27 = $p.Slides.Item(27)

# Duplicate the blank template slide twice: the first copy becomes the
# answer slide (28), the second stays blank as the new end-of-deck
# template (29).
$dup1 = $s27.Duplicate()
$s28 = $p.Slides.Item(28)
$dup2 = $s28.Duplicate()
$s29 = $p.Slides.Item(29)

# ---------------------------------------------------------------------
# 1. Slide 27 becomes the quiz question: drop the Title placeholder and
#    fill the Content placeholder with the question + four choices.
# ---------------------------------------------------------------------
$s27.Shapes.Item(1).Delete()

$tf27 = $s27.Shapes.Item(1).TextFrame
$tr27 = $tf27.TextRange
$tr27.Text = "Which of these is TRUE about Padding in "
$run = $tr27.InsertAfter("CNN?Padding")
$run = $run.InsertAfter(" is used in _________ layer.")
$run = $run.InsertAfter("`rConvolution as well as pooling")
$run = $run.InsertAfter("`rConvolution & Fully connected")
$run = $run.InsertAfter("`rFully connected & pooling")
$run = $run.InsertAfter("`rOnly convolution")
$run = $run.InsertAfter("`r")
$tf27.TextRange.Paragraphs(1, 1).Font.Bold = 1

# ---------------------------------------------------------------------
# 2. Slide 28 becomes the answer: drop the Title placeholder and fill
#    the Content placeholder with the answer sentence.
# ---------------------------------------------------------------------
$s28.Shapes.Item(1).Delete()

$tf28 = $s28.Shapes.Item(1).TextFrame
$tf28.TextRange.Text = "Padding is used in convolution as well as pooling layers in CNN."

# ---------------------------------------------------------------------
# 3. Slide 29 stays exactly as duplicated: a fresh, empty Title +
#    Content placeholder pair ready for the next quiz question.
# ---------------------------------------------------------------------
